$wb = $excel.ActiveWorkbook

# Update the "zh-cn" sheet's handoff/handback datetime stamps for row 2
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-25 10:39:39"
$wsZh.Range("H2").Value = "2016-03-25 10:40:18"

# Update the "de-de" sheet's handoff/handback datetime stamps for row 2
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-25 10:39:47"
$wsDe.Range("H2").Value = "2016-03-25 10:40:32"
